$wb = $excel.ActiveWorkbook

# Sheet1: "Heat Generators" - update energy carrier (input) column for rows 10-12 to "electricity"
$ws1 = $wb.Worksheets.Item("Heat Generators")
$ws1.Range("C10").Value = "electricity"
$ws1.Range("C11").Value = "electricity"
$ws1.Range("C12").Value = "electricity"
$ws1.Range("E31").Select()

# Sheet2: "prices and emmision factors" - set emission factor (C) for rows 6 and 7 to 1
$ws2 = $wb.Worksheets.Item("prices and emmision factors")
$ws2.Range("C6").Value = 1
$ws2.Range("C7").Value = 1
$ws2.Range("D24").Select()

# Sheet4: "Heat Storage" - update unloading efficiency (G3) from 0.52 to 0.95
$ws4 = $wb.Worksheets.Item("Heat Storage")
$ws4.Range("G3").Value = 0.95
$ws4.Range("I15").Select()
